# Append the "BASELINE" measurement rows (35-39) to the data_catalog sheet,
# plus a trailing blank row (40), matching the test-vector-gen commit for p3388.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values are written in the order that reproduces the author's shared
# string table growth (Excel interns strings in the order cells are
# committed, not sheet reading order): "none" first, then the five
# per-row "baseline_*.iq" file names (one per row, top to bottom), then
# "baseline measurement" and finally "BASELINE".

# Column D/F/G ("none") for all five rows
$ws.Range("D35").Value = "none"
$ws.Range("F35").Value = "none"
$ws.Range("G35").Value = "none"
$ws.Range("D36").Value = "none"
$ws.Range("F36").Value = "none"
$ws.Range("G36").Value = "none"
$ws.Range("D37").Value = "none"
$ws.Range("F37").Value = "none"
$ws.Range("G37").Value = "none"
$ws.Range("D38").Value = "none"
$ws.Range("F38").Value = "none"
$ws.Range("G38").Value = "none"
$ws.Range("D39").Value = "none"
$ws.Range("F39").Value = "none"
$ws.Range("G39").Value = "none"

# Column R / S (file names), one new string per row
$ws.Range("R35").Value = "baseline_2.4GHz_IQ_time.iq"
$ws.Range("S35").Value = "baseline_2.4GHz_IQ_time.iq"
$ws.Range("R36").Value = "baseline_5.3GHz_IQ_time.iq"
$ws.Range("S36").Value = "baseline_5.3GHz_IQ_time.iq"
$ws.Range("R37").Value = "baseline_0301_5GHz_IQ_time.iq"
$ws.Range("S37").Value = "baseline_0301_5GHz_IQ_time.iq"
$ws.Range("R38").Value = "baseline_0301_24GHz_IQ_time.iq"
$ws.Range("S38").Value = "baseline_0301_24GHz_IQ_time.iq"
$ws.Range("R39").Value = "baseline_0301_900_IQ_time.iq"
$ws.Range("S39").Value = "baseline_0301_900_IQ_time.iq"

# Column T ("baseline measurement") for all five rows
$ws.Range("T35").Value = "baseline measurement"
$ws.Range("T36").Value = "baseline measurement"
$ws.Range("T37").Value = "baseline measurement"
$ws.Range("T38").Value = "baseline measurement"
$ws.Range("T39").Value = "baseline measurement"

# Column A ("BASELINE") for all five rows
$ws.Range("A35").Value = "BASELINE"
$ws.Range("A36").Value = "BASELINE"
$ws.Range("A37").Value = "BASELINE"
$ws.Range("A38").Value = "BASELINE"
$ws.Range("A39").Value = "BASELINE"

# --- Remaining columns: numbers + already-interned strings (H/J/K/L/M). ---

$ws.Range("B35").Value = 100
$ws.Range("C35").Value = 1
$ws.Range("E35").Value = 0
$ws.Range("H35").Value = "OFF"
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = "no"
$ws.Range("K35").Value = "N/A"
$ws.Range("L35").Value = "N/A"
$ws.Range("M35").Value = "N/A"
$ws.Range("N35").Value = 2.5
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("Q35").Value = 4

$ws.Range("B36").Value = 100
$ws.Range("C36").Value = 2
$ws.Range("E36").Value = 0
$ws.Range("H36").Value = "OFF"
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = "no"
$ws.Range("K36").Value = "N/A"
$ws.Range("L36").Value = "N/A"
$ws.Range("M36").Value = "N/A"
$ws.Range("N36").Value = 5.3
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = 0
$ws.Range("Q36").Value = 4

$ws.Range("B37").Value = 100
$ws.Range("C37").Value = 3
$ws.Range("E37").Value = 0
$ws.Range("H37").Value = "OFF"
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = "no"
$ws.Range("K37").Value = "N/A"
$ws.Range("L37").Value = "N/A"
$ws.Range("M37").Value = "N/A"
$ws.Range("N37").Value = 5.3
$ws.Range("O37").Value = 0
$ws.Range("P37").Value = 0
$ws.Range("Q37").Value = 4

$ws.Range("B38").Value = 100
$ws.Range("C38").Value = 4
$ws.Range("E38").Value = 0
$ws.Range("H38").Value = "OFF"
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = "no"
$ws.Range("K38").Value = "N/A"
$ws.Range("L38").Value = "N/A"
$ws.Range("M38").Value = "N/A"
$ws.Range("N38").Value = 2.4
$ws.Range("O38").Value = 0
$ws.Range("P38").Value = 0
$ws.Range("Q38").Value = 4

$ws.Range("B39").Value = 100
$ws.Range("C39").Value = 5
$ws.Range("E39").Value = 0
$ws.Range("H39").Value = "OFF"
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = "no"
$ws.Range("K39").Value = "N/A"
$ws.Range("L39").Value = "N/A"
$ws.Range("M39").Value = "N/A"
$ws.Range("N39").Value = 0.9
$ws.Range("O39").Value = 0
$ws.Range("P39").Value = 0
$ws.Range("Q39").Value = 4

# --- Number formats (mirror the columns immediately above, rows 30-34). ---

$ws.Range("A35:A39").NumberFormat = "@"
$ws.Range("D35:D39").NumberFormat = "@"
$ws.Range("F35:G39").NumberFormat = "@"
$ws.Range("H35:H39").NumberFormat = "@"
$ws.Range("J35:J39").NumberFormat = "@"
$ws.Range("R35:R39").NumberFormat = "@"

$ws.Range("B35:C39").NumberFormat = "0"
$ws.Range("I35:I39").NumberFormat = "0"

$ws.Range("E35:E39").NumberFormat = "0.000"

$ws.Range("K35:M39").NumberFormat = "0"
$ws.Range("K35:M39").HorizontalAlignment = -4108  # xlCenter

$ws.Range("N35:Q39").NumberFormat = "0.0"

$ws.Range("T35").HorizontalAlignment = -4131  # xlLeft
$ws.Range("T36:T39").HorizontalAlignment = -4131  # xlLeft

# --- Trailing blank row 40 (A/B/C only, matching the author's row stub). ---
$ws.Range("A40").Value = ""
$ws.Range("B40").Value = ""
$ws.Range("C40").Value = ""
$ws.Range("A40").NumberFormat = "@"
$ws.Range("B40:C40").NumberFormat = "0"

# --- View bits that happen to be reproducible: the active selection. ---
$ws.Range("G29").Select()
